$d = $word.ActiveDocument

$table = $d.Tables.Item(1)

# Row 3 is the first empty data row following the header row and the
# "1 Change (New description)" / "PC_1" row.
$row = $table.Rows.Item(3)

$cell1 = $row.Cells.Item(1)
$range1 = $cell1.Range
$range1.Collapse(0)  # wdCollapseEnd -> position before the cell's end-of-cell mark
$range1.Text = "2 Change (Update)"
$range1.Font.Reset()
$range1.LanguageID = 1033

$cell2 = $row.Cells.Item(2)
$range2 = $cell2.Range
$range2.Collapse(0)  # wdCollapseEnd
$range2.Text = "PC_1"
$range2.Font.Reset()
$range2.LanguageID = 1033
